# Sprint 3 data entry — "Add files via upload / New Excel Sheet For Sprint 3"
#
# The workbook already contains a (blank) Sprint3 sheet/template; this
# commit fills in the actual Sprint-3 numbers (team man-hours table and the
# burndown task table) and makes Sprint3 the active/selected sheet, which is
# what a user does right after finishing data entry on that tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint3")

# --- Team man-hours table (rows 5-8: Cameron, Hassan, James, Tim) ---------
# Cameron
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 12

# Hassan
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 6

# James
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 6

# Tim — only week 3 hours changed (C8/D8 already populated)
$ws.Range("E8").Value = 32

# --- Burndown task table (rows 19-28): Remaining / Mon / Tue / Wed hours --
# Story 1: Priority Tasks
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 10

$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0

$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0

$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 2
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 4

$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 12

$ws.Range("G24").Value = 10
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0

$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 5

# Story 2: BurnDown/Backlog
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 1

# Story 3: Report
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2

$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 2

# --- Make Sprint3 the active sheet / selected cell, matching the author's
#     last on-screen state when the workbook was saved -------------------
$ws.Activate()
$ws.Range("H33").Select()
